$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended to the sheet (BOA sheet content additions)
$ws.Range("A12").Value = "Merrill Guided Investing4"

$ws.Range("A13").Value = "Review our Online Privacy Notice"
$ws.Range("A13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 17

$ws.Range("A14").Value = "Your security is our top priority"

# Slight column width adjustment (closest achievable quantized width)
$ws.Columns.Item(1).ColumnWidth = 29.83

# Move the active selection to the last entered cell, matching the saved view state
[void]$ws.Range("A14").Select()
